$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.851.02"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.639.30"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.20"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5073"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2583"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06444"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.66"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07789"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.290"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.865.58"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.633.34"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5634"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7608"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.23"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.883.54"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.23"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.335"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.902"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.122"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.774"
$ws.Range("E25").Value = "  -6.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1272"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "140.10"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.789"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.47"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04893"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.304"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.226"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.560"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.368"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9060"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.130.60"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5512"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.539"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8016"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.71"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.773.85"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  -6.44%  "
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.699"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05051"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.14%  "
